# Applies the updates described in the commit:
# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# (Remove previous account statements and add new ones, database modified)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Periodo Mora" column (E16:E22) with the new set of periods.
$ws.Range("E16").Value = "1910"
$ws.Range("E17").Value = "1909"
$ws.Range("E18").Value = "1908"
$ws.Range("E19").Value = "1907"
$ws.Range("E20").Value = "1906"
$ws.Range("E21").Value = "1904"
$ws.Range("E22").Value = "1807"

# Update "Valor Mora" column values that changed for the first/last rows.
$ws.Range("F16").Value = 42000
$ws.Range("F22").Value = 44000
